$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 1782.1
$ws.Range("I9").Value = 1782.1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1782.1
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -1613.1
# Row 10
$ws.Range("H10").Value = 10000
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -9707
$ws.Range("N10").Value = -10586
# Row 20
$ws.Range("H20").Value = 11007
$ws.Range("I20").Value = 2010.5
$ws.Range("J20").Value = 29000
$ws.Range("K20").Value = 2010.5
$ws.Range("L20").Value = 29000
$ws.Range("M20").Value = -1780.5
$ws.Range("N20").Value = -29460
# Row 34
$ws.Range("H34").Value = 4297.8237
$ws.Range("I34").Value = 1070.9333
$ws.Range("J34").Value = 28499.5
$ws.Range("K34").Value = 1070.9333
$ws.Range("L34").Value = 28499.5
$ws.Range("M34").Value = -867.9332999999999
$ws.Range("N34").Value = -28905.5
# Row 35
$ws.Range("H35").Value = 11007
$ws.Range("I35").Value = 2010.5
$ws.Range("J35").Value = 29000
$ws.Range("K35").Value = 2010.5
$ws.Range("L35").Value = 29000
$ws.Range("M35").Value = -1631.5
$ws.Range("N35").Value = -29758
# Row 36
$ws.Range("H36").Value = 4297.8237
$ws.Range("I36").Value = 1070.9333
$ws.Range("J36").Value = 28499.5
$ws.Range("K36").Value = 1070.9333
$ws.Range("L36").Value = 28499.5
$ws.Range("M36").Value = -355.9332999999999
$ws.Range("N36").Value = -29929.5
# Row 80
$ws.Range("H80").Value = 16393.154
$ws.Range("I80").Value = 6845.6665
$ws.Range("J80").Value = 37875
$ws.Range("K80").Value = 20536.9995
$ws.Range("L80").Value = 113625
$ws.Range("M80").Value = -19538.9995
$ws.Range("N80").Value = -115621
# Row 83
$ws.Range("H83").Value = 16393.154
$ws.Range("I83").Value = 6845.6665
$ws.Range("J83").Value = 37875
$ws.Range("K83").Value = 61610.9985
$ws.Range("L83").Value = 340875
$ws.Range("M83").Value = -56618.9985
$ws.Range("N83").Value = -350859
# Row 111
$ws.Range("H111").Value = 753
$ws.Range("I111").Value = 592.4286
$ws.Range("J111").Value = 865.4
$ws.Range("K111").Value = 1777.2858
$ws.Range("L111").Value = 2596.2
$ws.Range("M111").Value = 1289.7142
$ws.Range("N111").Value = -8730.2

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1101340.6
$ws.Range("I32").Value = 1369693.6
$ws.Range("J32").Value = 10038.267
$ws.Range("K32").Value = 1369693.6
$ws.Range("L32").Value = 10038.267
$ws.Range("M32").Value = -1369406.6
$ws.Range("N32").Value = -10612.267
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 18343.6
$ws.Range("I82").Value = 9860
$ws.Range("J82").Value = 26827.2
$ws.Range("K82").Value = 9860
$ws.Range("L82").Value = 26827.2
$ws.Range("M82").Value = -9477
$ws.Range("N82").Value = -27593.2
# Row 85
$ws.Range("H85").Value = 18343.6
$ws.Range("I85").Value = 9860
$ws.Range("J85").Value = 26827.2
$ws.Range("K85").Value = 9860
$ws.Range("L85").Value = 26827.2
$ws.Range("M85").Value = -8534
$ws.Range("N85").Value = -29479.2
# Row 126
$ws.Range("H126").Value = 24000
$ws.Range("J126").Value = 24000
$ws.Range("L126").Value = 24000
$ws.Range("N126").Value = -33880

$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 1400
$ws.Range("J14").Value = 2300
$ws.Range("K14").Value = 1400
$ws.Range("L14").Value = 2300
$ws.Range("M14").Value = -1230
$ws.Range("N14").Value = -2640
# Row 70
$ws.Range("H70").Value = 9000
$ws.Range("J70").Value = 9000
$ws.Range("L70").Value = 9000
$ws.Range("N70").Value = -9630
# Row 73
$ws.Range("H73").Value = 9000
$ws.Range("J73").Value = 9000
$ws.Range("L73").Value = 9000
$ws.Range("N73").Value = -11184
# Row 132
$ws.Range("H132").Value = 2397.1853
$ws.Range("I132").Value = 1583.25
$ws.Range("J132").Value = 3581.0908
$ws.Range("K132").Value = 4749.75
$ws.Range("L132").Value = 10743.2724
$ws.Range("M132").Value = -2219.75
$ws.Range("N132").Value = -15803.2724
# Row 134
$ws.Range("H134").Value = 1939.3334
$ws.Range("I134").Value = 958.64703
$ws.Range("J134").Value = 3221.7693
$ws.Range("K134").Value = 2875.94109
$ws.Range("L134").Value = 9665.3079
$ws.Range("M134").Value = -340.9410899999998
$ws.Range("N134").Value = -14735.3079

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 2223.4285
$ws.Range("I131").Value = 3130
$ws.Range("J131").Value = 1964.4082
$ws.Range("K131").Value = 9390
$ws.Range("L131").Value = 5893.2246
$ws.Range("M131").Value = -4350
$ws.Range("N131").Value = -15973.2246
# Row 132
$ws.Range("H132").Value = 6706.6665
$ws.Range("J132").Value = 16175
$ws.Range("L132").Value = 145575
$ws.Range("N132").Value = -150635

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 5775
$ws.Range("J15").Value = 5775
$ws.Range("L15").Value = 5775
$ws.Range("N15").Value = -6351
# Row 81
$ws.Range("H81").Value = 5775
$ws.Range("J81").Value = 5775
$ws.Range("L81").Value = 5775
$ws.Range("N81").Value = -7771
# Row 84
$ws.Range("H84").Value = 5775
$ws.Range("J84").Value = 5775
$ws.Range("L84").Value = 17325
$ws.Range("N84").Value = -27309

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1844.7646
$ws.Range("I16").Value = 1866.3125
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1866.3125
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1696.3125
$ws.Range("N16").Value = -1840
# Row 138
$ws.Range("H138").Value = 25000
$ws.Range("J138").Value = 25000
$ws.Range("L138").Value = 25000
$ws.Range("N138").Value = -35280

$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 10000
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11144
